$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '28.102.04'
Set-TextValue $ws 'E2' '  -1.82%  '
Set-TextValue $ws 'D3' '1.800.25'
Set-TextValue $ws 'E3' '  +0.13%  '
Set-TextValue $ws 'E4' '  +0.08%  '
Set-TextValue $ws 'D5' '316.67'
Set-TextValue $ws 'E5' '  +1.07%  '
Set-TextValue $ws 'E6' '  +0.03%  '
Set-TextValue $ws 'D7' '0.5451'
Set-TextValue $ws 'E7' '  +1.33%  '
Set-TextValue $ws 'D8' '0.3785'
Set-TextValue $ws 'E8' '  +0.21%  '
Set-TextValue $ws 'D9' '0.07463'
Set-TextValue $ws 'E9' '  -0.93%  '
Set-TextValue $ws 'D10' '41.99'
Set-TextValue $ws 'E10' '  -1.19%  '
Set-TextValue $ws 'E11' '  -2.17%  '
Set-TextValue $ws 'E12' '  +0.11%  '
Set-TextValue $ws 'D13' '6.194'
Set-TextValue $ws 'E13' '  +0.06%  '
Set-TextValue $ws 'E14' '  -2.62%  '
Set-TextValue $ws 'D15' '7.346'
Set-TextValue $ws 'E15' '  -1.26%  '
Set-TextValue $ws 'D16' '1.799.26'
Set-TextValue $ws 'E16' '  +0.14%  '
Set-TextValue $ws 'D17' '89.54'
Set-TextValue $ws 'E17' '  -0.81%  '
Set-TextValue $ws 'D18' '0.00001064'
Set-TextValue $ws 'E18' '  -0.05%  '
Set-TextValue $ws 'D19' '0.06544'
Set-TextValue $ws 'E19' '  +1.48%  '
Set-TextValue $ws 'D20' '17.41'
Set-TextValue $ws 'E20' '  +1.05%  '
Set-TextValue $ws 'E21' '  +0.05%  '
Set-TextValue $ws 'D22' '5.934'
Set-TextValue $ws 'E22' '  +0.03%  '
Set-TextValue $ws 'D23' '28.137.10'
Set-TextValue $ws 'E23' '  -1.72%  '
Set-TextValue $ws 'E24' '  +0.07%  '
Set-TextValue $ws 'D25' '2.090'
Set-TextValue $ws 'E25' '  -0.37%  '
Set-TextValue $ws 'D26' '155.61'
Set-TextValue $ws 'E26' '  -3.13%  '
Set-TextValue $ws 'D27' '20.42'
Set-TextValue $ws 'E27' '  -0.21%  '
Set-TextValue $ws 'D28' '2.005.72'
Set-TextValue $ws 'E28' '  -0.01%  '
Set-TextValue $ws 'D29' '2.322'
Set-TextValue $ws 'E29' '  -2.57%  '
Set-TextValue $ws 'D30' '121.80'
Set-TextValue $ws 'E30' '  -0.88%  '
Set-TextValue $ws 'D31' '0.1115'
Set-TextValue $ws 'E31' '  +8.54%  '
Set-TextValue $ws 'D32' '1.115'
Set-TextValue $ws 'E32' '  +0.73%  '
Set-TextValue $ws 'D33' '3.678'
Set-TextValue $ws 'E33' '  -0.12%  '
Set-TextValue $ws 'E34' '  -1.98%  '
Set-TextValue $ws 'D35' '0.06900'
Set-TextValue $ws 'E35' '  +6.74%  '
Set-TextValue $ws 'D36' '0.2222'
Set-TextValue $ws 'E36' '  -1.41%  '
Set-TextValue $ws 'D37' '0.02290'
Set-TextValue $ws 'E37' '  -0.48%  '
Set-TextValue $ws 'D38' '5.090'
Set-TextValue $ws 'E38' '  +0.92%  '
Set-TextValue $ws 'D39' '8.450'
Set-TextValue $ws 'E39' '  -4.91%  '
Set-TextValue $ws 'D40' '11.19'
Set-TextValue $ws 'E40' '  -1.48%  '
Set-TextValue $ws 'D41' '0.6157'
Set-TextValue $ws 'E41' '  -1.60%  '
Set-TextValue $ws 'B42' 'TrustWalletToken'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D42' '1.172'
Set-TextValue $ws 'E42' '  -3.17%  '
Set-TextValue $ws 'B43' 'WEMIXTOKEN'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D43' '1.419'
Set-TextValue $ws 'E43' '  +2.02%  '
Set-TextValue $ws 'D44' '13.29'
Set-TextValue $ws 'E44' '  -0.71%  '
Set-TextValue $ws 'D45' '3.682'
Set-TextValue $ws 'E45' '  +0.59%  '
Set-TextValue $ws 'D46' '0.5739'
Set-TextValue $ws 'E46' '  -2.54%  '
Set-TextValue $ws 'D47' '124.65'
Set-TextValue $ws 'E47' '  -1.08%  '
Set-TextValue $ws 'E48' '  +2.17%  '
Set-TextValue $ws 'D49' '1.919'
Set-TextValue $ws 'E49' '  -2.34%  '
Set-TextValue $ws 'D50' '0.06813'
Set-TextValue $ws 'E50' '  -1.32%  '
Set-TextValue $ws 'D51' '71.86'
Set-TextValue $ws 'E51' '  -1.09%  '
